# Insert a new data row at row 41 (pushing the existing row 41..115 down to
# 42..116), then populate the new row with the weekly entry that belongs
# between the current rows 40 (2021-02-09 / 44236) and the old row 41
# (2021-12-07 / 44537).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("41:41").Insert()

$ws.Range("A41").Value = 3
$ws.Range("B41").Value = "Femacal de La Calera"
$ws.Range("C41").Value = "Coquimbo"
$ws.Range("D41").Value = 44571
$ws.Range("E41").Value = 5
$ws.Range("F41").Value = 100112052
$ws.Range("G41").Value = "Albahaca"
$ws.Range("H41").Value = "Sin especificar"
$ws.Range("I41").Value = "Primera"
$ws.Range("J41").Value = 110
$ws.Range("K41").Value = 4500
$ws.Range("L41").Value = 5000
$ws.Range("M41").Value = 4773
$ws.Range("N41").Value = "`$/docena de matas"
$ws.Range("O41").Value = "Provincia de Quillota"
$ws.Range("P41").Value = 796
$ws.Range("Q41").Value = 6
$ws.Range("R41").Value = "Hortaliza"
